$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Seba"
$ws.Range("E3").Value = "Matias"
$ws.Range("F3").Value = "Javiera"

$ws.Range("F3").Select()
